$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Edit 1: merge "JEST" + " " + "Section 8 Course" runs into a single
# run reading " Section 8 Course" (same bold/underline/size formatting).
# -------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("JEST Section 8 Course", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $paraStart = $rng.Start
    $paraEnd = $rng.End
    $subStart = $paraStart + 4
    $subRng = $d.Range($subStart, $paraEnd)
    # First overwrite with a placeholder so the engine actually performs a
    # text replace (merging the run), then set the final text in a second
    # pass -- a no-op same-text assignment does not coalesce runs.
    $subRng.Text = "~~TEMP_MERGE_PLACEHOLDER~~"
    $subRng2 = $d.Range($subStart, $subStart + 27)
    $subRng2.Text = " Section 8 Course"
}

# -------------------------------------------------------------------
# Edit 2: append the new "Session Provider" notes section (Profile.tsx
# example + useSession explanation) right after the trailing code block,
# before the two blank paragraphs that close the document.
# -------------------------------------------------------------------
$count = $d.Paragraphs.Count
$spacer = $d.Paragraphs($count - 1)
$spacerRng = $spacer.Range
$insertPoint = $spacerRng.Duplicate
$insertPoint.Collapse(1)

$newXml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Session Provider:</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> (next-auth/react)</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">It is required when a component is client component and we </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>cant</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve"> use session like this directly </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">const session = await </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>auth(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>);</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">Since by default they are server components and we are using client component to fetch user’s session. </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t xml:space="preserve">In order to use it in server component: </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>We must wrap it using Session Provider, for example:</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t>const Profile = async () =&gt; {</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">  const session = await </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>auth(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>);</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">  return (</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    &lt;</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>SessionProvider</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> session={session}&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">      &lt;div </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>className</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>="max-w-md mx-auto space-y-4"&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        &lt;h2 </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>className</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>="h2-bold"&gt;Profile&lt;/h2&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">        &lt;</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ProfileForm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> /&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">      &lt;/div&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">    &lt;/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>SessionProvider</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">  );</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t>};</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t>export default Profile;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">The client component that is using this Session Provider for example in Profile Form component it must have </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
        <w:t>Use Session hook to consume it.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">import </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t xml:space="preserve">{ </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>useSession</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> } from "next-auth/react";</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:pStyle w:val="NoSpacing"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">  const </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>{ data</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve">: session, </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>update }</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> = </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>useSession</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>(</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t>);</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:sz w:val="28"/>
          <w:szCs w:val="28"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@

$insertPoint.InsertXML($newXml)
